$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.018.94"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.282.72"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'495.41"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "'127.71"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "2.281.57"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "'0.0949"
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "2.691.81"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "'21.74"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "54.129.96"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D18").Value = "2.375.01"
$ws.Range("E18").Value = "  +6.93%  "
$ws.Range("D19").Value = "'10.02"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("D20").Value = "'4.10"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").Value = "'301.03"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = "  +5.75%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'5.39"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "'62.30"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.396.38"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "'0.372"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'168.49"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'1.61"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "0.0₃0687"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").Value = "'5.87"
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "'17.71"
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").Value = "'0.905"
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").Value = "'3.69"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("D42").Value = "'35.55"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").Value = "'0.373"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").Value = "'3.36"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'126.98"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").Value = "'4.76"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("D48").Value = "'0.0887"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "'0.544"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'237.63"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").Value = "'0.0484"
$ws.Range("E51").Value = "  +3.46%  "
